$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.52139999999999
$ws.Range("D7").Value = -7.113100000000004
$ws.Range("A8").Value = -22.37970000000002
$ws.Range("A10").Value = -21.74160000000001
$ws.Range("A12").Value = -21.5788
$ws.Range("D15").Value = -7.894199999999999
$ws.Range("A18").Value = -22.38720000000002
$ws.Range("D18").Value = -8.214599999999997
$ws.Range("E18").Value = 16.4812
$ws.Range("E19").Value = 16.4964
$ws.Range("D20").Value = -7.863100000000002
$ws.Range("E27").Value = 16.53089999999999
$ws.Range("D29").Value = -6.8291
$ws.Range("D30").Value = -7.855799999999999
$ws.Range("D31").Value = -7.522399999999994
$ws.Range("E31").Value = 17.05710000000001
$ws.Range("A37").Value = -19.80599999999999
$ws.Range("E38").Value = 16.07119999999999
$ws.Range("D40").Value = -8.461199999999993
$ws.Range("E42").Value = 16.4387
$ws.Range("E44").Value = 16.5152
$ws.Range("E47").Value = 16.45639999999999
$ws.Range("D50").Value = -8.072499999999989
$ws.Range("A55").Value = -22.2678
$ws.Range("E58").Value = 16.32890000000001
$ws.Range("E65").Value = 17.285
$ws.Range("A68").Value = -21.5634
$ws.Range("D68").Value = -6.875799999999999
$ws.Range("E73").Value = 17.50770000000002
$ws.Range("D76").Value = -7.523599999999999
$ws.Range("A77").Value = -20.85969999999999
$ws.Range("A78").Value = -20.56069999999998
$ws.Range("A81").Value = -21.81910000000001
$ws.Range("A82").Value = -22.0729
$ws.Range("D87").Value = -7.944399999999994
$ws.Range("D88").Value = -7.283299999999997
$ws.Range("E90").Value = 16.30879999999999
$ws.Range("E94").Value = 19.05740000000003
$ws.Range("E95").Value = 18.23240000000002
$ws.Range("D96").Value = -7.6725
$ws.Range("D98").Value = -8.370100000000008
$ws.Range("D101").Value = -7.7802
$ws.Range("E101").Value = 16.668
$ws.Range("D102").Value = -8.015300000000002
